$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefix the "New nominations" / "Carryover nominations" / "Unconfirmed" /
# "Withdrawn" / "Confirmed" sub-rows in column A with the name of the
# section (category) they belong to, e.g. "     New nominations" becomes
# "     Civilian, New nominations".

$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Carryover nominations"
$ws.Range("A9").Value  = "     Civilian, Unconfirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn  "

$ws.Range("A12").Value = "     Civilian (lists), New nominations"
$ws.Range("A13").Value = "     Civilian (lists), Confirmed  "
$ws.Range("A14").Value = "     Civilian (lists), Unconfirmed  "

$ws.Range("A16").Value = "     Air Force, New nominations"
$ws.Range("A17").Value = "     Air Force, Carryover nominations"
$ws.Range("A18").Value = "     Air Force, Unconfirmed "
$ws.Range("A19").Value = "     Air Force, Withdrawn  "

$ws.Range("A21").Value = "     Army, New nominations"
$ws.Range("A22").Value = "     Army, Carryover nominations"
$ws.Range("A23").Value = "     Army, Unconfirmed "
$ws.Range("A24").Value = "     Army, Withdrawn  "

$ws.Range("A26").Value = "     Navy, New nominations"
$ws.Range("A27").Value = "     Navy, Carryover nominations"
$ws.Range("A28").Value = "     Navy, Unconfirmed  "

$ws.Range("A30").Value = "     Marine Corps, New nominations"
$ws.Range("A31").Value = "     Marine Corps, Carryover nominations"

# Rework the trailing "Summary" block: it used to be
#   32 Summary
#   33 Nominations carried over from first session    2762
#   34 Total nominations received this session         44934
#   35 Total confirmed                                  42493
#   36 Total unconfirmed                                 5189
#   37 Total withdrawn                                     14
# and becomes
#   32 Total new nominations                           44934
#   33 Total carryover nominations                       2762
#   34 Total confirmed                                   42493
#   35 Total unconfirmed                                  5189
#   36 Total withdrawn                                      14
# (one fewer row; values for confirmed/unconfirmed/withdrawn shift up).

$ws.Range("A32").Value = "Total new nominations"
# B32 is a brand-new number cell (row 32 used to be the bare "Summary"
# header with no value); pick up the "#,##0" look of its sibling total
# cells below (B33) before writing the value into it.
$ws.Range("B33").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$ws.Range("B32").Value = 44934

$ws.Range("A33").Value = "Total carryover nominations"
$ws.Range("B33").Value = 2762

$ws.Range("A34").Value = "Total confirmed  "
$ws.Range("B34").Value = 42493

$ws.Range("A35").Value = "Total unconfirmed  "
$ws.Range("B35").Value = 5189

$ws.Range("A36").Value = "Total withdrawn  "
# B36 carries over the "#,##0" format from the old "Total unconfirmed"
# row; the new "Total withdrawn" row should instead use the plain
# right-aligned look of the old row 37 (copied from B7 below) it replaces.
$ws.Range("B7").Copy()
$ws.Range("B36").PasteSpecial(-4122)
$ws.Range("B36").Value = 14

# The old row 37 ("Total withdrawn  " / 14) is now redundant - remove it.
$ws.Rows.Item(37).Delete()
